$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 161.292709
$ws.Range("H2").Value = 483.8781269999999
$ws.Range("I2").Value = 0.2350858997670666
$ws.Range("J2").Value = 0.2350858997670667
$ws.Range("M2").Value = 68.63737500000001
$ws.Range("N2").Value = 205.912125
$ws.Range("O2").Value = 0.5415701538216162
$ws.Range("P2").Value = 0.5415701538216162
$ws.Range("Q2").Value = 11070.70815239887
$ws.Range("R2").Value = 99636.37337158987
$ws.Range("S2").Value = 0.1273155068981433
$ws.Range("T2").Value = 0.1273155068981434
$ws.Range("G3").Value = 161.292709
$ws.Range("H3").Value = 483.8781269999999
$ws.Range("I3").Value = 0.2350858997670666
$ws.Range("J3").Value = 0.2350858997670667
$ws.Range("O3").Value = 0.08718851262838957
$ws.Range("P3").Value = 0.08718851262838957
$ws.Range("Q3").Value = 1782.296477638939
$ws.Range("R3").Value = 16040.66829875045
$ws.Range("S3").Value = 0.02049678994059722
$ws.Range("T3").Value = 0.02049678994059722
$ws.Range("G4").Value = 161.292709
$ws.Range("H4").Value = 483.8781269999999
$ws.Range("I4").Value = 0.2350858997670666
$ws.Range("J4").Value = 0.2350858997670667
$ws.Range("M4").Value = 16.21089566666667
$ws.Range("N4").Value = 48.632687
$ws.Range("O4").Value = 0.1279089892319285
$ws.Range("P4").Value = 0.1279089892319285
$ws.Range("Q4").Value = 2614.699277393028
$ws.Range("R4").Value = 23532.29349653725
$ws.Range("S4").Value = 0.03006959982188395
$ws.Range("T4").Value = 0.03006959982188396
$ws.Range("G5").Value = 161.292709
$ws.Range("H5").Value = 483.8781269999999
$ws.Range("I5").Value = 0.2350858997670666
$ws.Range("J5").Value = 0.2350858997670667
$ws.Range("M5").Value = 20.32546233333333
$ws.Range("N5").Value = 60.976387
$ws.Range("O5").Value = 0.1603741949973873
$ws.Range("P5").Value = 0.1603741949973873
$ws.Range("Q5").Value = 3278.348881420794
$ws.Range("R5").Value = 29505.13993278715
$ws.Range("S5").Value = 0.0377017119303798
$ws.Range("T5").Value = 0.03770171193037981
$ws.Range("G6").Value = 161.292709
$ws.Range("H6").Value = 483.8781269999999
$ws.Range("I6").Value = 0.2350858997670666
$ws.Range("J6").Value = 0.2350858997670667
$ws.Range("M6").Value = 10.513928
$ws.Range("N6").Value = 31.541784
$ws.Range("O6").Value = 0.08295814932067838
$ws.Range("P6").Value = 0.08295814932067838
$ws.Range("Q6").Value = 1695.819929350952
$ws.Range("R6").Value = 15262.37936415857
$ws.Range("S6").Value = 0.01950229117606235
$ws.Range("T6").Value = 0.01950229117606235
$ws.Range("I7").Value = 0.2000330076689186
$ws.Range("J7").Value = 0.2000330076689187
$ws.Range("M7").Value = 68.63737500000001
$ws.Range("N7").Value = 205.912125
$ws.Range("O7").Value = 0.5415701538216162
$ws.Range("P7").Value = 0.5415701538216162
$ws.Range("Q7").Value = 9419.990952002625
$ws.Range("R7").Value = 84779.91856802363
$ws.Range("S7").Value = 0.1083319067326568
$ws.Range("T7").Value = 0.1083319067326568
$ws.Range("I8").Value = 0.2000330076689186
$ws.Range("J8").Value = 0.2000330076689187
$ws.Range("O8").Value = 0.08718851262838957
$ws.Range("P8").Value = 0.08718851262838957
$ws.Range("S8").Value = 0.01744058041523626
$ws.Range("T8").Value = 0.01744058041523627
$ws.Range("I9").Value = 0.2000330076689186
$ws.Range("J9").Value = 0.2000330076689187
$ws.Range("M9").Value = 16.21089566666667
$ws.Range("N9").Value = 48.632687
$ws.Range("O9").Value = 0.1279089892319285
$ws.Range("P9").Value = 0.1279089892319285
$ws.Range("Q9").Value = 2224.829992461958
$ws.Range("R9").Value = 20023.46993215762
$ws.Range("S9").Value = 0.02558601982395399
$ws.Range("T9").Value = 0.02558601982395399
$ws.Range("I10").Value = 0.2000330076689186
$ws.Range("J10").Value = 0.2000330076689187
$ws.Range("M10").Value = 20.32546233333333
$ws.Range("N10").Value = 60.976387
$ws.Range("O10").Value = 0.1603741949973873
$ws.Range("P10").Value = 0.1603741949973873
$ws.Range("Q10").Value = 2789.524967632724
$ws.Range("R10").Value = 25105.72470869452
$ws.Range("S10").Value = 0.03208013257780903
$ws.Range("T10").Value = 0.03208013257780904
$ws.Range("I11").Value = 0.2000330076689186
$ws.Range("J11").Value = 0.2000330076689187
$ws.Range("M11").Value = 10.513928
$ws.Range("N11").Value = 31.541784
$ws.Range("O11").Value = 0.08295814932067838
$ws.Range("P11").Value = 0.08295814932067838
$ws.Range("Q11").Value = 1442.961748318712
$ws.Range("R11").Value = 12986.65573486841
$ws.Range("S11").Value = 0.01659436811926256
$ws.Range("T11").Value = 0.01659436811926256
$ws.Range("G12").Value = 160.3204953333334
$ws.Range("H12").Value = 480.961486
$ws.Range("I12").Value = 0.2336688876404109
$ws.Range("J12").Value = 0.2336688876404109
$ws.Range("M12").Value = 68.63737500000001
$ws.Range("N12").Value = 205.912125
$ws.Range("O12").Value = 0.5415701538216162
$ws.Range("P12").Value = 0.5415701538216162
$ws.Range("Q12").Value = 11003.97795837975
$ws.Range("R12").Value = 99035.80162541776
$ws.Range("S12").Value = 0.1265480954227433
$ws.Range("T12").Value = 0.1265480954227433
$ws.Range("G13").Value = 160.3204953333334
$ws.Range("H13").Value = 480.961486
$ws.Range("I13").Value = 0.2336688876404109
$ws.Range("J13").Value = 0.2336688876404109
$ws.Range("O13").Value = 0.08718851262838957
$ws.Range("P13").Value = 0.08718851262838957
$ws.Range("Q13").Value = 1771.553444030319
$ws.Range("R13").Value = 15943.98099627287
$ws.Range("S13").Value = 0.02037324276089771
$ws.Range("T13").Value = 0.02037324276089771
$ws.Range("G14").Value = 160.3204953333334
$ws.Range("H14").Value = 480.961486
$ws.Range("I14").Value = 0.2336688876404109
$ws.Range("J14").Value = 0.2336688876404109
$ws.Range("M14").Value = 16.21089566666667
$ws.Range("N14").Value = 48.632687
$ws.Range("O14").Value = 0.1279089892319285
$ws.Range("P14").Value = 0.1279089892319285
$ws.Range("Q14").Value = 2598.938823076988
$ws.Range("R14").Value = 23390.44940769289
$ws.Range("S14").Value = 0.02988835123303403
$ws.Range("T14").Value = 0.02988835123303403
$ws.Range("G15").Value = 160.3204953333334
$ws.Range("H15").Value = 480.961486
$ws.Range("I15").Value = 0.2336688876404109
$ws.Range("J15").Value = 0.2336688876404109
$ws.Range("M15").Value = 20.32546233333333
$ws.Range("N15").Value = 60.976387
$ws.Range("O15").Value = 0.1603741949973873
$ws.Range("P15").Value = 0.1603741949973873
$ws.Range("Q15").Value = 3258.58818915901
$ws.Range("R15").Value = 29327.29370243108
$ws.Range("S15").Value = 0.03747445975126585
$ws.Range("T15").Value = 0.03747445975126585
$ws.Range("G16").Value = 160.3204953333334
$ws.Range("H16").Value = 480.961486
$ws.Range("I16").Value = 0.2336688876404109
$ws.Range("J16").Value = 0.2336688876404109
$ws.Range("M16").Value = 10.513928
$ws.Range("N16").Value = 31.541784
$ws.Range("O16").Value = 0.08295814932067838
$ws.Range("P16").Value = 0.08295814932067838
$ws.Range("Q16").Value = 1685.598144859003
$ws.Range("R16").Value = 15170.38330373103
$ws.Range("S16").Value = 0.01938473847247002
$ws.Range("T16").Value = 0.01938473847247003
$ws.Range("G17").Value = 74.79809033333333
$ws.Range("H17").Value = 224.394271
$ws.Range("I17").Value = 0.1090190404506753
$ws.Range("J17").Value = 0.1090190404506753
$ws.Range("M17").Value = 68.63737500000001
$ws.Range("N17").Value = 205.912125
$ws.Range("O17").Value = 0.5415701538216162
$ws.Range("P17").Value = 0.5415701538216162
$ws.Range("Q17").Value = 5133.944575492876
$ws.Range("R17").Value = 46205.50117943588
$ws.Range("S17").Value = 0.05904145850635723
$ws.Range("T17").Value = 0.05904145850635724
$ws.Range("G18").Value = 74.79809033333333
$ws.Range("H18").Value = 224.394271
$ws.Range("I18").Value = 0.1090190404506753
$ws.Range("J18").Value = 0.1090190404506753
$ws.Range("O18").Value = 0.08718851262838957
$ws.Range("P18").Value = 0.08718851262838957
$ws.Range("Q18").Value = 826.5244831074116
$ws.Range("R18").Value = 7438.720347966704
$ws.Range("S18").Value = 0.009505207985068617
$ws.Range("T18").Value = 0.009505207985068619
$ws.Range("G19").Value = 74.79809033333333
$ws.Range("H19").Value = 224.394271
$ws.Range("I19").Value = 0.1090190404506753
$ws.Range("J19").Value = 0.1090190404506753
$ws.Range("M19").Value = 16.21089566666667
$ws.Range("N19").Value = 48.632687
$ws.Range("O19").Value = 0.1279089892319285
$ws.Range("P19").Value = 0.1279089892319285
$ws.Range("Q19").Value = 1212.544038459575
$ws.Range("R19").Value = 10912.89634613618
$ws.Range("S19").Value = 0.01394451527108061
$ws.Range("T19").Value = 0.01394451527108061
$ws.Range("G20").Value = 74.79809033333333
$ws.Range("H20").Value = 224.394271
$ws.Range("I20").Value = 0.1090190404506753
$ws.Range("J20").Value = 0.1090190404506753
$ws.Range("M20").Value = 20.32546233333333
$ws.Range("N20").Value = 60.976387
$ws.Range("O20").Value = 0.1603741949973873
$ws.Range("P20").Value = 0.1603741949973873
$ws.Range("Q20").Value = 1520.305767675431
$ws.Range("R20").Value = 13682.75190907888
$ws.Range("S20").Value = 0.01748384085166466
$ws.Range("T20").Value = 0.01748384085166466
$ws.Range("G21").Value = 74.79809033333333
$ws.Range("H21").Value = 224.394271
$ws.Range("I21").Value = 0.1090190404506753
$ws.Range("J21").Value = 0.1090190404506753
$ws.Range("M21").Value = 10.513928
$ws.Range("N21").Value = 31.541784
$ws.Range("O21").Value = 0.08295814932067838
$ws.Range("P21").Value = 0.08295814932067838
$ws.Range("Q21").Value = 786.4217363021627
$ws.Range("R21").Value = 7077.795626719464
$ws.Range("S21").Value = 0.009044017836504199
$ws.Range("T21").Value = 0.009044017836504201
$ws.Range("G22").Value = 152.4469883333333
$ws.Range("H22").Value = 457.340965
$ws.Range("I22").Value = 0.2221931644729284
$ws.Range("J22").Value = 0.2221931644729285
$ws.Range("M22").Value = 68.63737500000001
$ws.Range("N22").Value = 205.912125
$ws.Range("O22").Value = 0.5415701538216162
$ws.Range("P22").Value = 0.5415701538216162
$ws.Range("Q22").Value = 10463.56110585563
$ws.Range("R22").Value = 94172.04995270062
$ws.Range("S22").Value = 0.1203331862617155
$ws.Range("T22").Value = 0.1203331862617155
$ws.Range("G23").Value = 152.4469883333333
$ws.Range("H23").Value = 457.340965
$ws.Range("I23").Value = 0.2221931644729284
$ws.Range("J23").Value = 0.2221931644729285
$ws.Range("O23").Value = 0.08718851262838957
$ws.Range("P23").Value = 0.08718851262838957
$ws.Range("Q23").Value = 1684.550603791796
$ws.Range("R23").Value = 15160.95543412616
$ws.Range("S23").Value = 0.01937269152658976
$ws.Range("T23").Value = 0.01937269152658977
$ws.Range("G24").Value = 152.4469883333333
$ws.Range("H24").Value = 457.340965
$ws.Range("I24").Value = 0.2221931644729284
$ws.Range("J24").Value = 0.2221931644729285
$ws.Range("M24").Value = 16.21089566666667
$ws.Range("N24").Value = 48.632687
$ws.Range("O24").Value = 0.1279089892319285
$ws.Range("P24").Value = 0.1279089892319285
$ws.Range("Q24").Value = 2471.302222569218
$ws.Range("R24").Value = 22241.72000312295
$ws.Range("S24").Value = 0.02842050308197592
$ws.Range("T24").Value = 0.02842050308197593
$ws.Range("G25").Value = 152.4469883333333
$ws.Range("H25").Value = 457.340965
$ws.Range("I25").Value = 0.2221931644729284
$ws.Range("J25").Value = 0.2221931644729285
$ws.Range("M25").Value = 20.32546233333333
$ws.Range("N25").Value = 60.976387
$ws.Range("O25").Value = 0.1603741949973873
$ws.Range("P25").Value = 0.1603741949973873
$ws.Range("Q25").Value = 3098.555519199273
$ws.Range("R25").Value = 27886.99967279346
$ws.Range("S25").Value = 0.03563404988626798
$ws.Range("T25").Value = 0.03563404988626799
$ws.Range("G26").Value = 152.4469883333333
$ws.Range("H26").Value = 457.340965
$ws.Range("I26").Value = 0.2221931644729284
$ws.Range("J26").Value = 0.2221931644729285
$ws.Range("M26").Value = 10.513928
$ws.Range("N26").Value = 31.541784
$ws.Range("O26").Value = 0.08295814932067838
$ws.Range("P26").Value = 0.08295814932067838
$ws.Range("Q26").Value = 1602.816659153507
$ws.Range("R26").Value = 14425.34993238156
$ws.Range("S26").Value = 0.01843273371637925
$ws.Range("T26").Value = 0.01843273371637925
